# Daily attendance processing - 2026-01-30 12:08:03
# Swap the order of the "Recorded By" names in column G wherever a session
# was recorded by both the instructor account and the System (i.e. cells
# whose text is exactly "dnasr281@gmail.com, System"), turning them into
# "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
